{"js": "const replacements = [\n  [\"2025-10-03 Friday\", \"2025-10-04 Saturday\"],\n  [\"325\u00d79=2925\", \"549\u00d75=2745\"],\n  [\"208\u00d78=1664\", \"694\u00d77=4858\"],\n  [\"867\u00d76=5202\", \"217\u00d79=1953\"],\n  [\"967\u00d74=3868\", \"249\u00d77=1743\"],\n  [\"545\u00d72=1090\", \"511\u00d78=4088\"],\n  [\"158\u00d76=948\", \"478\u00d76=2868\"],\n  [\"726\u00d77=5082\", \"386\u00d73=1158\"],\n  [\"276\u00d74=1104\", \"280\u00d72=560\"],\n  [\"392\u00d74=1568\", \"703\u00d74=2812\"],\n  [\"268\u00d75=1340\", \"922\u00d74=3688\"],\n  [\"677\u00d77=4739\", \"113\u00d79=1017\"],\n  [\"602\u00d75=3010\", \"197\u00d75=985\"],\n  [\"142\u00d75=710\", \"422\u00d74=1688\"],\n  [\"652\u00d72=1304\", \"949\u00d77=6643\"],\n  [\"815\u00d74=3260\", \"264\u00d78=2112\"],\n  [\"517\u00d74=2068\", \"783\u00d76=4698\"],\n  [\"535\u00d72=1070\", \"381\u00d79=3429\"],\n  [\"494\u00d73=1482\", \"648\u00d73=1944\"],\n  [\"143\u00d78=1144\", \"779\u00d76=4674\"],\n  [\"750\u00d76=4500\", \"812\u00d76=4872\"],\n  [\"346\u00d72=692\", \"813\u00d76=4878\"],\n  [\"536\u00d73=1608\", \"270\u00d76=1620\"],\n  [\"380\u00d74=1520\", \"395\u00d75=1975\"],\n  [\"345\u00d77=2415\", \"420\u00d73=1260\"],\n  [\"516\u00d73=1548\", \"336\u00d73=1008\"],\n];\n\nfor (const [find, replace] of replacements) {\n  const results = context.document.body.search(find, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length === 0) {\n    throw new Error(\"No match found for: \" + find);\n  }\n\n  for (const item of results.items) {\n    item.insertText(replace, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "$d = $word.ActiveDocument\n\n$replacements = @(\n    @(\"2025-10-03 Friday\", \"2025-10-04 Saturday\"),\n    @(\"325\u00d79=2925\", \"549\u00d75=2745\"),\n    @(\"208\u00d78=1664\", \"694\u00d77=4858\"),\n    @(\"867\u00d76=5202\", \"217\u00d79=1953\"),\n    @(\"967\u00d74=3868\", \"249\u00d77=1743\"),\n    @(\"545\u00d72=1090\", \"511\u00d78=4088\"),\n    @(\"158\u00d76=948\", \"478\u00d76=2868\"),\n    @(\"726\u00d77=5082\", \"386\u00d73=1158\"),\n    @(\"276\u00d74=1104\", \"280\u00d72=560\"),\n    @(\"392\u00d74=1568\", \"703\u00d74=2812\"),\n    @(\"268\u00d75=1340\", \"922\u00d74=3688\"),\n    @(\"677\u00d77=4739\", \"113\u00d79=1017\"),\n    @(\"602\u00d75=3010\", \"197\u00d75=985\"),\n    @(\"142\u00d75=710\", \"422\u00d74=1688\"),\n    @(\"652\u00d72=1304\", \"949\u00d77=6643\"),\n    @(\"815\u00d74=3260\", \"264\u00d78=2112\"),\n    @(\"517\u00d74=2068\", \"783\u00d76=4698\"),\n    @(\"535\u00d72=1070\", \"381\u00d79=3429\"),\n    @(\"494\u00d73=1482\", \"648\u00d73=1944\"),\n    @(\"143\u00d78=1144\", \"779\u00d76=4674\"),\n    @(\"750\u00d76=4500\", \"812\u00d76=4872\"),\n    @(\"346\u00d72=692\", \"813\u00d76=4878\"),\n    @(\"536\u00d73=1608\", \"270\u00d76=1620\"),\n    @(\"380\u00d74=1520\", \"395\u00d75=1975\"),\n    @(\"345\u00d77=2415\", \"420\u00d73=1260\"),\n    @(\"516\u00d73=1548\", \"336\u00d73=1008\"),\n)\n\nforeach ($pair in $replacements) {\n    $findText = $pair[0]\n    $replaceText = $pair[1]\n\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $findText\n    $find.Replacement.Text = $replaceText\n\n    $result = $find.Execute(\n        $findText, $false, $false, $false, $false, $false, $true, 1, $true, $replaceText, 2\n    )\n\n    if (-not $result) {\n        Write-Output \"WARNING: replacement failed for $findText\"\n    }\n}\n"}
